# Re-order the "Recorded By" (column G) comma-separated author lists so
# that the last author in the list becomes the first, preserving the
# relative order of the remaining authors. Cells with only a single
# value (no comma) are left untouched. Row 1 is the header ("Recorded By")
# and is skipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G is the 7th column ("Recorded By").
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($r -eq 1) {
        continue
    }

    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    if ($val -notmatch ",") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Length -le 1) {
        continue
    }

    $lastPart = $parts[$parts.Length - 1]
    $restParts = $parts[0..($parts.Length - 2)]
    $newParts = @($lastPart) + $restParts
    $newVal = [string]::Join(", ", $newParts)

    $cell.Value2 = $newVal
}
